$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 1911
$ws.Range("L3").Value = 1938
$ws.Range("L4").Value = 540
$ws.Range("L6").Value = 1763
$ws.Range("L7").Value = 6265

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 108
$ws.Range("L3").Value = 132
$ws.Range("L7").Value = 390

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 68
$ws.Range("L3").Value = 97
$ws.Range("L6").Value = 96
$ws.Range("L7").Value = 280

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L6").Value = 70
$ws.Range("L7").Value = 225

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 48
$ws.Range("L4").Value = 21
$ws.Range("L8").Value = 390
$ws.Range("L11").Value = 115
$ws.Range("L16").Value = 13
$ws.Range("L19").Value = 177
$ws.Range("L20").Value = 165
$ws.Range("L23").Value = 66
$ws.Range("L27").Value = 65
$ws.Range("L29").Value = 323
$ws.Range("L33").Value = 280
$ws.Range("L36").Value = 90
$ws.Range("L37").Value = 225
$ws.Range("L48").Value = 92
$ws.Range("L50").Value = 38
$ws.Range("L52").Value = 131
$ws.Range("L54").Value = 134
$ws.Range("L55").Value = 57
$ws.Range("L63").Value = 17
$ws.Range("L67").Value = 221
$ws.Range("L71").Value = 17
$ws.Range("L72").Value = 27
$ws.Range("L74").Value = 5
$ws.Range("L77").Value = 40
$ws.Range("L79").Value = 176
$ws.Range("L80").Value = 23
$ws.Range("L85").Value = 332
$ws.Range("L86").Value = 48
$ws.Range("L88").Value = 82
$ws.Range("L89").Value = 76
$ws.Range("L93").Value = 33
$ws.Range("L98").Value = 47
$ws.Range("L101").Value = 6265

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 70
$ws.Range("L6").Value = 60
$ws.Range("L7").Value = 221

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L4").Value = 11
$ws.Range("L7").Value = 134

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 102
$ws.Range("L3").Value = 115
$ws.Range("L4").Value = 12
$ws.Range("L6").Value = 88
$ws.Range("L7").Value = 323

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 20
$ws.Range("L7").Value = 92

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L3").Value = 54
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 177

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 49
$ws.Range("L4").Value = 20

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 23
$ws.Range("L3").Value = 20
$ws.Range("L7").Value = 57

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 21
$ws.Range("L7").Value = 66

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 63
$ws.Range("L3").Value = 59
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 176

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 55
$ws.Range("L3").Value = 47
$ws.Range("L7").Value = 165

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L2").Value = 11
$ws.Range("L3").Value = 8
$ws.Range("L7").Value = 33

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L3").Value = 8
$ws.Range("L7").Value = 47

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 38

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 115

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L2").Value = 21
$ws.Range("L7").Value = 82

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 28
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 76

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 27
$ws.Range("L6").Value = 48

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 137
$ws.Range("L7").Value = 332

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L2").Value = 4
$ws.Range("L7").Value = 17

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 27

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 40

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L3").Value = 6
$ws.Range("L6").Value = 23

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 44
$ws.Range("L7").Value = 131

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L3").Value = 8
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L6").Value = 10
$ws.Range("L7").Value = 13

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 5

Write-Output "All 2025-04-27 updates applied."